$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting the existing row 86 (and all
# rows below it) down by one. Excel re-bases all the following rows, which
# also pushes the former last data row (126) into a brand-new row 127.
$ws.Rows(86).Insert()

# Populate the newly inserted row 86 with the new weekly record.
$ws.Range("A86").Value = 10
$ws.Range("B86").Value = "Vega Modelo de Temuco"
$ws.Range("C86").Value = "La Araucanía"
$ws.Range("D86").Value = 44460
$ws.Range("E86").Value = 9
$ws.Range("F86").Value = 100112005
$ws.Range("G86").Value = "Puerro"
$ws.Range("H86").Value = "Azul de Maquehue"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 20
$ws.Range("K86").Value = 8000
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = 8000
$ws.Range("N86").Value = "$/docena de paquetes"
$ws.Range("O86").Value = "Provincia de Cautín"
$ws.Range("P86").Value = 667
$ws.Range("Q86").Value = 12
$ws.Range("R86").Value = "Hortaliza"
